$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Row 5 - assignment page / IHC QC / Informal Consults test
$ws.Range("A5").Value = "Test to see is the assignment page comes up on IHC QC and Informal Consults"
$ws.Range("B5").Value = "Create a MM case with an IHC, and one with an Informal Consult and proceed to where the assignment page should be."
$ws.Range("C5").Value = "One worked one did not"
$ws.Range("D5").Value = "Sid was notified and said he would fix the problem."

# Row 6 - Bladder TUR template test
$ws.Range("A6").Value = "Test the added Bladder TUR template to make sure that the Block submission code works properly."
$ws.Range("B6").Value = "create a MM case with a Bladder TUR template and multiple blocks"
$ws.Range("C6").Value = "Everything worked correctly"
$ws.Range("D6").Value = "N/A"

# Row 7 - linking / special characters test (B7 entered before A7 to match shared-string order)
$ws.Range("B7").Value = "create multiple MM cases with each of the special characters in question and one or two with no spaces or anything like is normal now and see if the linking page finds all of the patient cases to be linked."
$ws.Range("A7").Value = "Check to see how linking works with hyphens(-), apostrophies('), and spaces in the patient name and with multiple letters in the middle initial space"
$ws.Range("C7").Value = "works as Sid expected it to"
$ws.Range("D7").Value = "N/A"

# Match row heights to the new content
$ws.Rows.Item(5).RowHeight = 90
$ws.Rows.Item(6).RowHeight = 60
$ws.Rows.Item(7).RowHeight = 120

# Make "Random Testing" the active sheet/tab and select E7
$ws.Select()
$ws.Range("E7").Select()
